# Add 2022-Q3 data
# 1) Insert a new worksheet (copy of "2022-Q2") positioned right before "2022-Q2",
#    rename it to "2022-Q3" and replace its contents with the new quarter's fund
#    holdings data.
# 2) Update the "总计" (summary) sheet: shift the existing quarter rows down by one
#    and insert the new 2022-Q3 summary row at the top of the data (row 2), adding
#    one extra row at the bottom for the quarter that got pushed out.

$wb = $excel.ActiveWorkbook

$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet)

$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# Clear out the copied 2022-Q2 data rows, keep header (row 1) as-is.
$newSheet.Range("A2:H3").ClearContents()

# Make sure every data row (2-12) has the same formatting as a normal data row:
# column A bold/bordered like the header, columns B-H unformatted.
$newSheet.Range("A2:H2").Copy() | Out-Null
$newSheet.Range("A2:H12").PasteSpecial(-4122) | Out-Null

$fundRows = @(
    @(0,  "290011", "泰信中小盘精选混合",     "11.26", "94.48", "9.15", "1.0303", 7),
    @(1,  "270021", "广发聚瑞混合A",          "17.58", "93.91", "4.04", "0.7102", 9),
    @(2,  "506007", "广发科创板两年定开混合", "5.01",  "94.25", "4.40", "0.2204", 9),
    @(3,  "012342", "广发瑞泽精选混合A",      "4.99",  "93.90", "3.75", "0.1871", 9),
    @(4,  "002580", "泰信鑫选灵活配置混合C",  "1.06",  "94.19", "9.45", "0.1002", 7),
    @(5,  "013000", "广发盛泽一年持有期混合A", "2.29", "82.39", "4.13", "0.0946", 7),
    @(6,  "001970", "泰信鑫选灵活配置混合A",  "0.79",  "94.19", "9.45", "0.0747", 7),
    @(7,  "002133", "广发鑫益灵活配置混合",   "1.22",  "93.85", "4.97", "0.0606", 6),
    @(8,  "010026", "广发聚瑞混合C",          "0.50",  "93.91", "4.04", "0.0202", 9),
    @(9,  "012343", "广发瑞泽精选混合C",      "0.36",  "93.90", "3.75", "0.0135", 9),
    @(10, "013001", "广发盛泽一年持有期混合C", "0.27", "82.39", "4.13", "0.0112", 7)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[1]

    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[3]

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[4]

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[5]

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $row[6]

    $newSheet.Cells.Item($r, 8).Value = $row[7]

    # columns B-G hold text values; drop the forced "@" number format so the
    # saved style matches a plain/unstyled cell while keeping the text type.
    $newSheet.Range($newSheet.Cells.Item($r, 2), $newSheet.Cells.Item($r, 7)).Style = "Normal"

    $r = $r + 1
}

$newSheet.Range("A1").Select() | Out-Null

# --- Update the 总计 (summary) sheet -------------------------------------

$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @("2022-Q3", 11, 2.52),
    @("2022-Q2", 2, 0.37),
    @("2022-Q1", 12, 2.31),
    @("2021-Q4", 13, 4.45),
    @("2021-Q3", 11, 4.35),
    @("2021-Q2", 1, 0.01),
    @("2020-Q4", 1, 1.06)
)

# Extend formatting (column A bold/bordered) down to the new row 8 before
# writing values, by copying the existing row 7 formatting down one row.
$summary.Range("A7:D7").Copy() | Out-Null
$summary.Range("A8:D8").PasteSpecial(-4122) | Out-Null

$r = 2
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# Row 8's index cell (A8) continues the 0-based sequence started in A2.
$summary.Cells.Item(8, 1).Value = 6

$summary.Range("A1").Select() | Out-Null
